$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.512.27"
$ws.Range("E2").Value = "  +8.02%  "
$ws.Range("D3").Value = "3.408.95"
$ws.Range("E3").Value = "  +4.95%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "411.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.84%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "123.20"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +14.34%  "
$ws.Range("D7").Value = "3.402.48"
$ws.Range("E7").Value = "  +4.88%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.582"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.45%  "
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.641"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.85%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.123"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +29.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "41.16"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.10%  "
$ws.Range("E13").Value = "  -0.66%  "
$ws.Range("D14").Value = "3.946.05"
$ws.Range("E14").Value = "  +4.98%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.41"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.68%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.52"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.49%  "
$ws.Range("D17").Value = "3.409.22"
$ws.Range("E17").Value = "  +4.78%  "
$ws.Range("D18").Value = "61.421.23"
$ws.Range("E18").Value = "  +8.18%  "
$ws.Range("E19").Value = "  -0.65%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.83"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0000123"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +14.01%  "
$ws.Range("E22").Value = "  -0.42%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.82"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.31%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "298.92"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.74%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "76.58"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.45%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.14"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "30.66"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +9.59%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.05"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +10.41%  "
$ws.Range("E29").Value = "  -2.32%  "
$ws.Range("E30").Value = "  -5.44%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.172"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.98%  "
$ws.Range("E32").Value = "  +6.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "42.73"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.11%  "
$ws.Range("B34").Value = "Toncoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.55"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +19.78%  "
$ws.Range("B35").Value = "Cosmos"
$ws.Range("C35").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.39"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0478"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "52.39"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.48%  "
$ws.Range("E39").Value = "  +1.69%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.01"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.63%  "
$ws.Range("E42").Value = "  +5.64%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.123"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.29%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "134.35"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.73%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "17.31"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.23%  "
$ws.Range("E46").Value = "  -0.60%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.91"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.08%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "21.85"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.46%  "
$ws.Range("B49").Value = "WEMIXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.18"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.33%  "
$ws.Range("D50").Value = "2.197.56"
$ws.Range("E50").Value = "  +2.47%  "
$ws.Range("D51").Value = "3.745.51"
$ws.Range("E51").Value = "  +5.03%  "
